# The author finalized ~10 years of BMKG daily-weather data: the raw data
# table (A9:K39, header + 30 daily rows) on the original "Data Harian -
# Table" sheet was selected, copied, and pasted (values + formats) onto a
# brand-new "Sheet1" tab (added right after the original), starting at A1.
# The new sheet becomes the active/selected tab.

$wb  = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item(1)

# Keep gridlines visible on the source sheet's view (matches the workbook's
# original display settings).
$excel.ActiveWindow.DisplayGridlines = $true

# Select + copy the finalized data range (header row 9 + data rows 10-39).
$srcRange = $src.Range("A9:K39")
$srcRange.Select()
$srcRange.Copy()

# Add the destination sheet right after the source sheet.
$dst = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $src)

# Paste values first ...
$dst.Range("A1").PasteSpecial()

# ... then copy the formats (borders / alignment / wrap) across separately
# so the header keeps its centered style and the data rows keep their
# left/top wrapped style, same as the source.
$src.Range("A9:K9").Copy()
$dst.Range("A1:K1").PasteSpecial(-4122)

$src.Range("A10:K39").Copy()
$dst.Range("A2:K31").PasteSpecial(-4122)

# Mirror the final selection left on each sheet.
$dst.Range("A1:K31").Select()
$dst.Activate()
